# The commit swaps the two embedded themes in the deck: the theme that is
# actually bound to the (single) slide master -- "Integral" -- is replaced
# by the stock "Office Theme" color palette (fonts/format scheme were
# already identical between the two themes, so only the 12 theme colors
# actually change in substance).
#
# PowerPoint's object model doesn't expose a "rewrite this theme part"
# verb, so we drive it the way a user would from the Design > Variants >
# Colors > Customize Colors dialog: push the 12 standard "Office" RGB
# values into the theme color scheme that the slide master (and therefore
# every slide) resolves to.

function Pack-OleColor([int]$r, [int]$g, [int]$b) {
    # PowerPoint/VBA RGB()-style packed color: 0x00BBGGRR
    return ($b * 65536) + ($g * 256) + $r
}

function Set-ThemeColor($themeColorScheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $themeColorScheme.Colors($index).RGB = Pack-OleColor $r $g $b
}

$p = $ppt.ActivePresentation

# Any slide's ThemeColorScheme resolves back to the shared slide-master
# theme (theme1.xml) -- there is only one design/master in this deck.
$slide = $p.Slides.Item(1)
$theme = $slide.ThemeColorScheme

# Standard Office theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeColor $theme 1  "000000"
Set-ThemeColor $theme 2  "FFFFFF"
Set-ThemeColor $theme 3  "44546A"
Set-ThemeColor $theme 4  "E7E6E6"
Set-ThemeColor $theme 5  "5B9BD5"
Set-ThemeColor $theme 6  "ED7D31"
Set-ThemeColor $theme 7  "A5A5A5"
Set-ThemeColor $theme 8  "FFC000"
Set-ThemeColor $theme 9  "4472C4"
Set-ThemeColor $theme 10 "70AD47"
Set-ThemeColor $theme 11 "0563C1"
Set-ThemeColor $theme 12 "954F72"
